# repull data, push all data, mean calculation
# Update the dSF column (F) values with refreshed/repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 9
    3  = 5
    4  = 5
    6  = 3
    7  = 1
    8  = 1
    9  = -4
    10 = 13
    11 = -1
    12 = 5
    13 = -1
    14 = 2
    16 = -5
    17 = 2
    18 = 5
    19 = -4
    21 = -4
    22 = 4
    23 = -6
    24 = 1
    25 = 6
    27 = 3
    28 = -1
    29 = -3
    31 = 4
    33 = 5
    34 = 1
    35 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
